# Update latest output (run 176)

$wb = $excel.ActiveWorkbook

# ---- Sheet "Schedule" ----
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E3").Value = 752.5666394999998
$schedule.Range("F3").Value = 28.44167193877551
$schedule.Range("E4").Value = 48.96743474999998
$schedule.Range("F4").Value = 1.439371979717813

# ---- Sheet "Detailed" ----
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B43").Value = 61.94424
$detailed.Range("B44").Value = 57.04367
$detailed.Range("B45").Value = 51.47522
$detailed.Range("C45").Value = "historical"
$detailed.Range("B46").Value = 50.55368
$detailed.Range("C46").Value = "historical"
$detailed.Range("B47").Value = 57.04365
$detailed.Range("B51").Value = 56.98
$detailed.Range("B54").Value = 36.2
$detailed.Range("B55").Value = 56.11627
$detailed.Range("B56").Value = 56.87432
$detailed.Range("B58").Value = 56.98
$detailed.Range("B61").Value = 57.06
$detailed.Range("B62").Value = 57.68364
$detailed.Range("B64").Value = 26.93884
$detailed.Range("B65").Value = 24.77568
$detailed.Range("B66").Value = 5.52969
$detailed.Range("B68").Value = 0.04948
$detailed.Range("B69").Value = 0.02903
$detailed.Range("B70").Value = 0.51
$detailed.Range("B71").Value = 0.51
$detailed.Range("B72").Value = 0.7
$detailed.Range("B73").Value = 0.7
$detailed.Range("B74").Value = 0.7
$detailed.Range("B75").Value = 0.02712
$detailed.Range("B76").Value = 0.02662
$detailed.Range("B77").Value = -1.30554
$detailed.Range("B78").Value = 0.008529999999999999
$detailed.Range("B79").Value = -2.54301
$detailed.Range("B80").Value = -5.01
$detailed.Range("B81").Value = -2.12343
$detailed.Range("B82").Value = -4.38184
$detailed.Range("B85").Value = -4.03608
$detailed.Range("B86").Value = 12.13262
$detailed.Range("B87").Value = 32.28228
$detailed.Range("B88").Value = 47.41492
$detailed.Range("B92").Value = 47.52905
$detailed.Range("B93").Value = 54.69705
$detailed.Range("B94").Value = 50.99335
$detailed.Range("B97").Value = 52.10975
